$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.927.08"
Set-TextValue $ws.Range("E2") "  +3.73%  "

Set-TextValue $ws.Range("D3") "3.101.69"
Set-TextValue $ws.Range("E3") "  +2.14%  "

Set-TextValue $ws.Range("E4") "  +0.10%  "

Set-TextValue $ws.Range("D5") "561.49"
Set-TextValue $ws.Range("E5") "  +3.07%  "

Set-TextValue $ws.Range("D6") "144.49"
Set-TextValue $ws.Range("E6") "  +7.77%  "

Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.04%  "

Set-TextValue $ws.Range("D8") "3.100.66"
Set-TextValue $ws.Range("E8") "  +2.34%  "

Set-TextValue $ws.Range("D9") "0.501"
Set-TextValue $ws.Range("E9") "  +1.30%  "

Set-TextValue $ws.Range("D10") "6.37"
Set-TextValue $ws.Range("E10") "  +3.48%  "

Set-TextValue $ws.Range("D11") "0.153"
Set-TextValue $ws.Range("E11") "  +3.70%  "

Set-TextValue $ws.Range("D12") "0.475"
Set-TextValue $ws.Range("E12") "  +6.33%  "

Set-TextValue $ws.Range("D13") "0.0000230"
Set-TextValue $ws.Range("E13") "  +3.36%  "

Set-TextValue $ws.Range("D14") "35.46"
Set-TextValue $ws.Range("E14") "  +3.70%  "

Set-TextValue $ws.Range("D15") "3.602.03"
Set-TextValue $ws.Range("E15") "  +2.37%  "

Set-TextValue $ws.Range("D16") "64.979.82"
Set-TextValue $ws.Range("E16") "  +3.91%  "

Set-TextValue $ws.Range("D17") "3.098.14"
Set-TextValue $ws.Range("E17") "  +2.59%  "

Set-TextValue $ws.Range("D18") "0.110"
Set-TextValue $ws.Range("E18") "  +1.70%  "

Set-TextValue $ws.Range("D19") "6.80"
Set-TextValue $ws.Range("E19") "  +2.50%  "

Set-TextValue $ws.Range("D20") "480.93"
Set-TextValue $ws.Range("E20") "  +0.74%  "

Set-TextValue $ws.Range("D21") "13.84"
Set-TextValue $ws.Range("E21") "  +4.35%  "

Set-TextValue $ws.Range("D22") "0.689"
Set-TextValue $ws.Range("E22") "  +2.32%  "

Set-TextValue $ws.Range("D23") "7.60"
Set-TextValue $ws.Range("E23") "  +7.53%  "

Set-TextValue $ws.Range("D24") "13.55"
Set-TextValue $ws.Range("E24") "  +11.96%  "

Set-TextValue $ws.Range("D25") "81.23"
Set-TextValue $ws.Range("E25") "  +0.12%  "

Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  -0.06%  "

Set-TextValue $ws.Range("D27") "2.80"
Set-TextValue $ws.Range("E27") "  +3.54%  "

Set-TextValue $ws.Range("D28") "8.21"
Set-TextValue $ws.Range("E28") "  +5.13%  "

Set-TextValue $ws.Range("D29") "2.07"
Set-TextValue $ws.Range("E29") "  +7.10%  "

Set-TextValue $ws.Range("D30") "0.997"
Set-TextValue $ws.Range("E30") "  -0.13%  "

Set-TextValue $ws.Range("D31") "26.21"
Set-TextValue $ws.Range("E31") "  +1.88%  "

Set-TextValue $ws.Range("E32") "  +1.91%  "

Set-TextValue $ws.Range("D33") "2.51"
Set-TextValue $ws.Range("E33") "  +5.75%  "

Set-TextValue $ws.Range("D34") "5.65"
Set-TextValue $ws.Range("E34") "  +0.03%  "

Set-TextValue $ws.Range("D35") "6.18"
Set-TextValue $ws.Range("E35") "  +5.30%  "

Set-TextValue $ws.Range("D36") "54.96"
Set-TextValue $ws.Range("E36") "  -0.03%  "

Set-TextValue $ws.Range("D37") "472.84"
Set-TextValue $ws.Range("E37") "  +2.27%  "

Set-TextValue $ws.Range("D38") "0.0839"
Set-TextValue $ws.Range("E38") "  +4.50%  "

Set-TextValue $ws.Range("D39") "0.0409"
Set-TextValue $ws.Range("E39") "  +5.33%  "

Set-TextValue $ws.Range("D40") "2.96"
Set-TextValue $ws.Range("E40") "  +20.24%  "

Set-TextValue $ws.Range("D41") "2.987.13"
Set-TextValue $ws.Range("E41") "  -5.52%  "

Set-TextValue $ws.Range("D42") "8.27"
Set-TextValue $ws.Range("E42") "  +2.21%  "

Set-TextValue $ws.Range("E43") "  -2.26%  "

Set-TextValue $ws.Range("D44") "28.41"
Set-TextValue $ws.Range("E44") "  +6.81%  "

Set-TextValue $ws.Range("D45") "0.261"
Set-TextValue $ws.Range("E45") "  +6.48%  "

Set-TextValue $ws.Range("D46") "2.17"
Set-TextValue $ws.Range("E46") "  +8.76%  "

Set-TextValue $ws.Range("E47") "  +0.00%  "

Set-TextValue $ws.Range("D48") "0.113"
Set-TextValue $ws.Range("E48") "  +3.68%  "

Set-TextValue $ws.Range("D49") "0.0₃0534"
Set-TextValue $ws.Range("E49") "  +6.26%  "

Set-TextValue $ws.Range("D50") "117.57"
Set-TextValue $ws.Range("E50") "  +2.76%  "

Set-TextValue $ws.Range("D51") "2.08"
Set-TextValue $ws.Range("E51") "  +3.67%  "
